$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Term Type" values in column C for rows 3-6 (keep the style)
$ws.Range("C3:C6").ClearContents()

# Delete the now-empty trailing rows 7-9
$ws.Range("A7:R9").Delete()

# Update the active selection
$ws.Range("C5").Select() | Out-Null
